# Weekly crime-data refresh for 107th Precinct CompStat report.
# Moves the reporting week forward by one (8/12-8/18/2024 -> 8/19-8/25/2024,
# report "Number" 33 -> 34) and updates the weekly/28-day/YTD/2-year/14-year/
# 31-year crime-count statistics and their derived percent-change figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings (report number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 16.666666666666
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = 55.555555555555
$ws.Range("N15").Value = -36.363636363636

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -61.111111111111
$ws.Range("I16").Value = 76
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = -18.279569892473
$ws.Range("L16").Value = -39.682539682539
$ws.Range("M16").Value = -54.761904761904
$ws.Range("N16").Value = -87.183811129848

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 21.052631578947
$ws.Range("I17").Value = 187
$ws.Range("J17").Value = 146
$ws.Range("K17").Value = 28.082191780821
$ws.Range("L17").Value = 53.27868852459
$ws.Range("M17").Value = 187.692307692308
$ws.Range("N17").Value = 16.875

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -31.578947368421
$ws.Range("I18").Value = 139
$ws.Range("J18").Value = 179
$ws.Range("K18").Value = -22.346368715083
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -22.346368715083
$ws.Range("N18").Value = -85.873983739837

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -17.647058823529
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -16.39344262295
$ws.Range("I19").Value = 367
$ws.Range("J19").Value = 398
$ws.Range("K19").Value = -7.788944723618
$ws.Range("L19").Value = -13.647058823529
$ws.Range("M19").Value = 25.255972696245
$ws.Range("N19").Value = -6.377551020408

# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 41
$ws.Range("H20").Value = -34.146341463414
$ws.Range("I20").Value = 229
$ws.Range("J20").Value = 245
$ws.Range("K20").Value = -6.530612244897
$ws.Range("L20").Value = 45.859872611465
$ws.Range("M20").Value = 69.629629629629
$ws.Range("N20").Value = -92.734771573604

# Row 21
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 49
$ws.Range("E21").Value = -24.489795918367
$ws.Range("F21").Value = 123
$ws.Range("G21").Value = 160
$ws.Range("H21").Value = -23.125
$ws.Range("I21").Value = 1015
$ws.Range("J21").Value = 1073
$ws.Range("K21").Value = -5.405405405405
$ws.Range("L21").Value = 3.045685279187
$ws.Range("M21").Value = 19.131455399061
$ws.Range("N21").Value = -80.87431693989

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -85.714285714285
$ws.Range("I22").Value = 25
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -13.793103448275
$ws.Range("L22").Value = 38.888888888888
$ws.Range("M22").Value = 56.25

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 38
$ws.Range("J23").Value = 54
$ws.Range("K23").Value = -29.629629629629
$ws.Range("L23").Value = -5
$ws.Range("M23").Value = 65.217391304347

# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -2.752293577981
$ws.Range("I24").Value = 783
$ws.Range("J24").Value = 922
$ws.Range("K24").Value = -15.075921908893
$ws.Range("L24").Value = -15.896885069817
$ws.Range("M24").Value = 20.461538461538

# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 331
$ws.Range("J25").Value = 330
$ws.Range("K25").Value = 0.30303030303
$ws.Range("L25").Value = 16.140350877193

# Row 26
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 5.263157894736
$ws.Range("I26").Value = 333
$ws.Range("J26").Value = 305
$ws.Range("K26").Value = 9.180327868852
$ws.Range("L26").Value = 11.744966442953
$ws.Range("M26").Value = 21.532846715328

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 21.052631578947
$ws.Range("L27").Value = -4.166666666666

# Row 28
$ws.Range("D28").Value = 3
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 6
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = -12.903225806451
$ws.Range("L28").Value = -32.5

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G31").Value = 1
$ws.Range("G31").NumberFormat = "#,##0"
$ws.Range("H31").Value = -100
$ws.Range("H31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J31").Value = 19
$ws.Range("K31").Value = -73.684210526315

# Row 33
$ws.Range("C33").Value = 1
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("F33").Value = 1
$ws.Range("F33").NumberFormat = "#,##0"
$ws.Range("I33").Value = 3
$ws.Range("K33").Value = 50
$ws.Range("L33").Value = -25

